# Voeg een nieuwe rij toe aan de namenlijst: Senna, ja (karten), nee (eigen vervoer)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Senna"
$ws.Range("B6").Value = "ja"
$ws.Range("C6").Value = "nee"

# Selectie verplaatsen zoals in het origineel (gebruiker klikte elders verder)
$ws.Range("K5").Select() | Out-Null
